$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "tasa"
$ws.Range("B2").Select()
